$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Add the new rows (8-14) of settings that describe the "beslutsdatum Fran"
# validation/logic outcomes. Values are entered in the same left-to-right /
# top-to-bottom order the author used (row 9 got its Value column filled in
# before its Name column, and the long comment in C8 was typed in last),
# so new shared-string entries line up with how the workbook was produced.
$ws.Range("A8").Value = "CaseNotValidYet"
$ws.Range("B8").Value = "Ärendet är inte giltigt ännu"

$ws.Range("B9").Value = "Ny eller växelvis adress, men elev + VH finns redan i filen"
$ws.Range("A9").Value = "NewOrVäxelvisAdress"

$ws.Range("A10").Value = "StudentAndHomeAdressAlreadyExsistInFile"
$ws.Range("B10").Value = "Nyare ärende med Elev + Adress finns redan i filen"

$ws.Range("A11").Value = "StudentAppearsTwoTimes"
$ws.Range("B11").Value = "Eleven har två giltiga ärenden i filen"

$ws.Range("A12").Value = "CaseNoLongerValid"
$ws.Range("B12").Value = "Beslutsdatumet är inte längre giltigt"

$ws.Range("A13").Value = "CaseIsNOTApprovedDecisionText"
$ws.Range("B13").Value = "Ärendet blev ej beviljat"

$ws.Range("A14").Value = "CaseIsNoLongerValidTooOld"
$ws.Range("B14").Value = "Ärendet är för gammalt"

$ws.Range("C8").Value = "Bör ej ändras, ärenden med denna kommentar sorteras ut för att behandlas när ärendet är nära sitt besluts datum"

# Columns got noticeably wider once the new, longer text was added -
# re-fit them to the new content.
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()

# The author ended up on the Settings sheet looking at the newly added
# row (A8) instead of the Assets sheet - move the active tab / selection
# there to match.
[void]$ws.Range("A8").Select()
